$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: find the hyperlink Address (target URL) bound to a given cell
# address (e.g. "$A$2") on a worksheet.
# ---------------------------------------------------------------------------
function Get-HLAddress($ws, $cellAddr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $cellAddr) {
            return $hl.Address
        }
    }
    return $null
}

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: the Status column for both rows (zh-cn + de-de columns)
# moves from "Ready for handoff" to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# Per-language detail sheets (zh-cn, de-de): the handback report now fills
# in "Latest Target File" (E) / "Latest Handback File" (F) with hyperlinks
# mirroring the source (A) / handoff (C) file links, refreshes the
# "Latest Handback DateTime" (G) and updates the shared Status text (B).
# ---------------------------------------------------------------------------
$langSheets = @("zh-cn", "de-de")
$handbackTimes = @{ "zh-cn" = "2016-01-18 04:42:54"; "de-de" = "2016-01-18 04:43:18" }

foreach ($langName in $langSheets) {
    $ws = $wb.Worksheets.Item($langName)

    $urlA2 = Get-HLAddress $ws "`$A`$2"
    $urlC2 = Get-HLAddress $ws "`$C`$2"
    $urlA3 = Get-HLAddress $ws "`$A`$3"
    $urlC3 = Get-HLAddress $ws "`$C`$3"

    $nameA2 = $ws.Range("A2").Value()
    $nameC2 = $ws.Range("C2").Value()
    $nameA3 = $ws.Range("A3").Value()
    $nameC3 = $ws.Range("C3").Value()

    $handbackTime = $handbackTimes[$langName]

    # Row 2
    $ws.Range("B2").Value = $newStatus
    $ws.Range("E2").Value = $nameA2
    $ws.Range("F2").Value = $nameC2
    $ws.Range("G2").Value = $handbackTime
    $ws.Hyperlinks.Add($ws.Range("E2"), $urlA2, [type]::Missing, [type]::Missing, $nameA2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $urlC2, [type]::Missing, [type]::Missing, $nameC2) | Out-Null
    $ws.Range("E2").Style = "HyperLink"
    $ws.Range("F2").Style = "HyperLink"

    # Row 3
    $ws.Range("B3").Value = $newStatus
    $ws.Range("E3").Value = $nameA3
    $ws.Range("F3").Value = $nameC3
    $ws.Range("G3").Value = $handbackTime
    $ws.Hyperlinks.Add($ws.Range("E3"), $urlA3, [type]::Missing, [type]::Missing, $nameA3) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $urlC3, [type]::Missing, [type]::Missing, $nameC3) | Out-Null
    $ws.Range("E3").Style = "HyperLink"
    $ws.Range("F3").Style = "HyperLink"
}

Write-Host "Report generated for handback"
